# Update crypto price/volume data per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain numeric-looking strings (e.g. "1.00").
# Force those specific cells to text format (one at a time -- applying
# NumberFormat to a multi-area union range does not reliably stick to
# every area) so Excel keeps the value as a string, matching the sheet's
# original inline-string cell type, instead of coercing it to a number.

$ws.Range('D2').Value = '60.237.16'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '2.408.17'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '558.94'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '135.64'
$ws.Range('E6').Value = '  -1.59%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.61'
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '24.72'
$ws.Range('E13').Value = '  -2.79%  '
$ws.Range('D14').Value = '2.837.30'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').Value = '60.136.88'
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('E16').Value = '  +0.37%  '
$ws.Range('D17').Value = '2.416.33'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.21'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.55'
$ws.Range('E19').Value = '  +3.19%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '325.58'
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.83'
$ws.Range('E21').Value = '  +1.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '64.60'
$ws.Range('E23').Value = '  -2.53%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.174'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.55'
$ws.Range('E25').Value = '  -2.46%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.80'
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('D29').Value = '0.0₃0771'
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.12'
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('E32').Value = '  +5.48%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.403'
$ws.Range('E33').Value = '  -2.00%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '18.39'
$ws.Range('E34').Value = '  -1.54%  '
$ws.Range('E35').Value = '  +3.62%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.17'
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '324.40'
$ws.Range('E39').Value = '  +3.42%  '
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '38.61'
$ws.Range('E41').Value = '  -2.58%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '148.71'
$ws.Range('E42').Value = '  +6.58%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.56'
$ws.Range('E43').Value = '  -3.39%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0970'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '19.89'
$ws.Range('E45').Value = '  +2.00%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0516'
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.575'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0221'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '11.05'
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('E50').Value = '  -1.04%  '
$ws.Range('E51').Value = '  -0.64%  '

# Restore the default (unstyled) cell style on the text-forced cells so
# no stray explicit number format lingers on them, matching the original
# workbook formatting.
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"

